# Apply the "cicada_ig" regeneration edits to ValueSet-interval-reason.xlsx
#
# Summary of the change (per the commit's xml diff):
#  1. Rename sheet 2 ("Include from Interval Reason" -> "Include #0")
#  2. On the Metadata sheet:
#       - URL value: pythia -> cicada
#       - Date value: updated to the newly generated timestamp
#       - A new "Jurisdiction" row (with an empty value) is inserted
#         right after the "Contact" row, shifting Description/Purpose/
#         Copyright/Immutable down by one row
#  3. On the Include sheet:
#       - The CodeSystem "System URI" value: pythia -> cicada
#         (all other rows on that sheet stay the same)

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item(1)
$wsInclude = $wb.Worksheets.Item(2)

# --- 1. Rename the "Include" worksheet tab -----------------------------
$wsInclude.Name = "Include #0"

# --- 2. Update the Metadata sheet --------------------------------------

# URL (row 2, column B)
$wsMeta.Range("B2").Value2 = "http://fhirfli.dev/fhir/ig/cicada/ValueSet/interval-reason"

# Date (row 8, column B)
$wsMeta.Range("B8").Value2 = "2026-02-11T14:37:07-05:00"

# Insert a new blank row at row 11 (before "Description"), shifting the
# rows below (Description, Purpose, Copyright, Immutable) down by one.
# -4121 = xlShiftDown
$wsMeta.Range("A11:B11").Insert(-4121)

# Copy the formatting from the row above (Contact, row 10) onto the new
# row so the inserted cells keep the same style as the rest of the table.
$wsMeta.Range("A10:B10").Copy()
$wsMeta.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new "Jurisdiction" property with an empty value.
$wsMeta.Range("A11").Value2 = "Jurisdiction"
$wsMeta.Range("B11").Value2 = ""

# --- 3. Update the Include sheet ---------------------------------------

# System URI (row 6, column B)
$wsInclude.Range("B6").Value2 = "http://fhirfli.dev/fhir/ig/cicada/CodeSystem/IntervalReason"

Write-Output "done"
